# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value for "展览" sheet (sheet1)
$exhibitUpdates = @{
    4  = 13498
    6  = 1025
    7  = 17
    8  = 1735
    9  = 135
    11 = 79
    14 = 13508
    15 = 336
    16 = 600
    17 = 8953
    19 = 8033
    20 = 251
    22 = 148
    27 = 1021
    31 = 205
    32 = 175
    34 = 95
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new value for "全部类型" sheet (sheet4)
$allUpdates = @{
    4  = 13498
    6  = 1025
    7  = 17
    8  = 1735
    9  = 135
    11 = 79
    14 = 13508
    15 = 336
    16 = 600
    17 = 8953
    19 = 8033
    22 = 148
    27 = 1021
    33 = 205
    34 = 175
    36 = 95
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
